$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (style) of the last existing year row (A10) onto the
# new row label cell A11 before we set its value.
$ws.Range("A10").Copy()
$ws.Range("A11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A11").Value = "2021年"

$rowData = @{
    "B11"  = 1687.64
    "C11"  = 488.36
    "D11"  = 84.86
    "F11"  = 1051.59
    "G11"  = 1742.13
    "H11"  = 188.41
    "I11"  = 946.76
    "J11"  = 322.88
    "K11"  = 373.95
    "L11"  = 139.12
    "M11"  = 22.97
    "N11"  = 476.78
    "O11"  = 720.58
    "P11"  = 65.88
    "Q11"  = 305.09
    "R11"  = 1146.82
    "S11"  = 30.34
    "T11"  = 1175.05
    "U11"  = 4.53
    "V11"  = 288.47
    "W11"  = 64.16
    "X11"  = 150.75
    "Y11"  = 2712.44
    "Z11"  = 369.38
    "AA11" = 290.88
    "AB11" = 1.04
    "AC11" = 27875.03
    "AD11" = 872.9400000000001
    "AE11" = 522.61
    "AF11" = 2830.79
    "AG11" = 1907.93
    "AH11" = 377.32
    "AI11" = 334.4
    "AJ11" = 28.34
    "AK11" = 1643.23
    "AL11" = 375.44
    "AM11" = 2213.19
    "AN11" = 174.15
    "AO11" = 546.04
    "AP11" = 1080.61
    "AQ11" = 116.97
}

foreach ($addr in $rowData.Keys) {
    $ws.Range($addr).Value = $rowData[$addr]
}

# E11 is an explicit, empty TEXT cell in the target row -- not simply an
# absent cell (assigning "" via .Value clears/removes the cell entirely,
# which is not what we want here). Using a bare quote-prefix formula
# forces an empty text entry to be recorded, then resetting the style
# back to Normal drops the "quote prefix" formatting flag it introduces.
$ws.Range("E11").Formula = "'"
$ws.Range("E11").Style = "Normal"
